$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17, shifting existing rows 17:67 down to 18:68.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new record.
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 44620
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 100112021
$ws.Range("G17").Value = "Ají"
$ws.Range("H17").Value = "Americana (o)"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 9000
$ws.Range("L17").Value = 9000
$ws.Range("M17").Value = 9000
$ws.Range("N17").Value = "$/caja 15 kilos"
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 600
$ws.Range("Q17").Value = 15
$ws.Range("R17").Value = "Hortaliza"
